$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("B2").Value = 46587.3723067849
$ws.Range("E2").Value = 51196.1892930378
$ws.Range("F2").Value = 53708.3527843201
$ws.Range("I2").Value = 6916.3723067849
$ws.Range("B3").Value = 40880.0398238708
$ws.Range("C3").Value = 35398.2597798799
$ws.Range("D3").Value = 33753.3437526218
$ws.Range("I3").Value = 19432.0398238708
$ws.Range("B4").Value = 54961.0256676177
$ws.Range("C4").Value = 48019.6394282117
$ws.Range("D4").Value = 46041.6316924572
$ws.Range("I4").Value = 13807.0256676177
$ws.Range("B5").Value = 53208.8556693793
$ws.Range("C5").Value = 47627.8956849143
$ws.Range("D5").Value = 45702.5029661378
$ws.Range("I5").Value = 6480.85566937931
$ws.Range("B6").Value = 55380.1398370414
$ws.Range("C6").Value = 49221.533962318
$ws.Range("D6").Value = 46762.9630901841
$ws.Range("I6").Value = 8627.13983704143
$ws.Range("B7").Value = 54132.1549147528
$ws.Range("C7").Value = 47326.3808502299
$ws.Range("D7").Value = 45042.5042915925
$ws.Range("I7").Value = 7594.15491475281
$ws.Range("B8").Value = 57788.9304815289
$ws.Range("C8").Value = 51550.2617191523
$ws.Range("D8").Value = 49648.3170404454
$ws.Range("I8").Value = 10950.9304815289
$ws.Range("B9").Value = 57336.3867463526
$ws.Range("C9").Value = 51043.3584047773
$ws.Range("D9").Value = 48617.7999808013
$ws.Range("I9").Value = 10498.3867463526
$ws.Range("B10").Value = 55811.7721773807
$ws.Range("C10").Value = 50993.9086660262
$ws.Range("D10").Value = 48528.0208235751
$ws.Range("I10").Value = 6846.77217738073
$ws.Range("B11").Value = 53490.5662733249
$ws.Range("E11").Value = 59806.6295152657
$ws.Range("F11").Value = 62607.8348206132
$ws.Range("I11").Value = 9052.56627332485
$ws.Range("B12").Value = 53831.9317275381
$ws.Range("I12").Value = 8526.93172753807
$ws.Range("B13").Value = 53116.2908240343
$ws.Range("I13").Value = 8420.29082403432
$ws.Range("B14").Value = 50375.1737232324
$ws.Range("E14").Value = 57033.6725959699
$ws.Range("F14").Value = 60406.8508024676
$ws.Range("I14").Value = 10296.1737232324
$ws.Range("B15").Value = 44877.0415480318
$ws.Range("C15").Value = 39493.7625990074
$ws.Range("D15").Value = 36620.7555123896
$ws.Range("E15").Value = 50538.1119257871
$ws.Range("F15").Value = 53076.3055806049
$ws.Range("I15").Value = 11262.0415480318
$ws.Range("B16").Value = 58919.2553935056
$ws.Range("C16").Value = 49638.4110313454
$ws.Range("D16").Value = 46608.2517880631
$ws.Range("I16").Value = 8237.25539350563
$ws.Range("B17").Value = 56864.8590570676
$ws.Range("C17").Value = 48512.9752276283
$ws.Range("D17").Value = 45926.8012163
$ws.Range("I17").Value = 7751.85905706759
$ws.Range("B18").Value = 59123.4642480049
$ws.Range("C18").Value = 50097.5266750737
$ws.Range("D18").Value = 48019.1181431368
$ws.Range("I18").Value = 11124.4642480049
$ws.Range("B19").Value = 57865.9146749818
$ws.Range("C19").Value = 49130.6985454978
$ws.Range("D19").Value = 46975.9889889691
$ws.Range("I19").Value = 10442.9146749818
$ws.Range("B20").Value = 61241.0936174398
$ws.Range("C20").Value = 52184.0019221491
$ws.Range("D20").Value = 49391.0500648189
$ws.Range("I20").Value = 9710.09361743976
$ws.Range("B21").Value = 60759.6737059826
$ws.Range("C21").Value = 51787.4509278864
$ws.Range("D21").Value = 49464.8868088423
$ws.Range("I21").Value = 14668.6737059826
$ws.Range("B22").Value = 59334.1815641534
$ws.Range("C22").Value = 51770.1963604418
$ws.Range("D22").Value = 49524.175585873
$ws.Range("I22").Value = 13542.1815641534
$ws.Range("B23").Value = 56981.3524785629
$ws.Range("I23").Value = 16081.3524785629
$ws.Range("B24").Value = 57393.7477228114
$ws.Range("I24").Value = 15219.7477228114
$ws.Range("B25").Value = 56727.6052436585
$ws.Range("I25").Value = 14406.6052436585
$ws.Range("B26").Value = 54078.3799783448
$ws.Range("I26").Value = 14492.3799783448
$ws.Range("B27").Value = 49036.3748764608
$ws.Range("I27").Value = 14353.3748764608
$ws.Range("B28").Value = 62578.6187859992
$ws.Range("C28").Value = 51629.3006580522
$ws.Range("D28").Value = 47959.7473417696
$ws.Range("I28").Value = 15600.6187859992
$ws.Range("B29").Value = 60342.3392906831
$ws.Range("C29").Value = 49693.4992795452
$ws.Range("D29").Value = 46890.4675071195
$ws.Range("I29").Value = 20829.3392906831
$ws.Range("B30").Value = 62677.9537925539
$ws.Range("C30").Value = 51865.9066446398
$ws.Range("D30").Value = 48928.8678881641
$ws.Range("I30").Value = 18926.9537925539
$ws.Range("B31").Value = 61381.4170032993
$ws.Range("C31").Value = 50966.4554825506
$ws.Range("D31").Value = 48302.1117238534
$ws.Range("I31").Value = 12874.4170032993
$ws.Range("B32").Value = 64593.458062834
$ws.Range("C32").Value = 52900.4000173944
$ws.Range("D32").Value = 50222.9814298912
$ws.Range("I32").Value = 13202.458062834
$ws.Range("B33").Value = 64075.3928887673
$ws.Range("C33").Value = 52482.2806734703
$ws.Range("D33").Value = 49734.6384747389
$ws.Range("I33").Value = 13593.3928887673
$ws.Range("B34").Value = 62648.126492201
$ws.Range("C34").Value = 52671.7780907208
$ws.Range("I34").Value = 18178.126492201
$ws.Range("B35").Value = 60184.6026306049
$ws.Range("I35").Value = 21130.6026306049
$ws.Range("B36").Value = 60675.5428037335
$ws.Range("I36").Value = 25523.5428037335
$ws.Range("B37").Value = 60031.1910642567
$ws.Range("I37").Value = 35664.1910642567
$ws.Range("B38").Value = 57364.020579266
$ws.Range("I38").Value = 28656.020579266
$ws.Range("B39").Value = 52752.3143768248
$ws.Range("I39").Value = 9178.31437682476
$ws.Range("B40").Value = 65834.2323540432
$ws.Range("C40").Value = 52324.9025529785
$ws.Range("I40").Value = 15979.2323540432
$ws.Range("B41").Value = 63602.059682333
$ws.Range("C41").Value = 51647.036593035
$ws.Range("I41").Value = 14676.059682333
$ws.Range("B42").Value = 65976.4073627165
$ws.Range("C42").Value = 52870.6715846101
$ws.Range("I42").Value = 12718.4073627165
$ws.Range("B43").Value = 64658.089745557
$ws.Range("C43").Value = 52308.4488349807
$ws.Range("I43").Value = 12651.089745557
$ws.Range("B44").Value = 67838.3330456273
$ws.Range("C44").Value = 53060.7347435344
$ws.Range("I44").Value = 9591.33304562733
$ws.Range("B45").Value = 67298.3204880051
$ws.Range("C45").Value = 53128.2230178537
$ws.Range("I45").Value = 6230.32048800515
$ws.Range("B46").Value = 65823.2641768403
$ws.Range("I46").Value = 10056.2641768403
$ws.Range("B47").Value = 63257.3249657043
$ws.Range("I47").Value = 6276.32496570428
$ws.Range("B48").Value = 63795.5800737923
$ws.Range("I48").Value = 6076.58007379228
$ws.Range("B49").Value = 63135.8908141192
$ws.Range("I49").Value = 12312.8908141192
